$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Distribution"
$ws.Range("A7").Font.Name = "Calibri"
$ws.Range("A7").Font.Size = 11
$ws.Range("A7").Font.Color = 0
